$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be stored as literal text (prevents Excel from
    # auto-converting numeric-looking strings into numbers/dates), then
    # reset the cell style back to Normal so no stray style index is left
    # behind (matches cells that carry no explicit style in the source).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '66.183.61'
Set-TextValue $ws.Range("E2") '  -1.24%  '
Set-TextValue $ws.Range("D3") '3.074.24'
Set-TextValue $ws.Range("E3") '  -1.64%  '
Set-TextValue $ws.Range("E4") '  +0.05%  '
Set-TextValue $ws.Range("D5") '573.94'
Set-TextValue $ws.Range("E5") '  -0.64%  '
Set-TextValue $ws.Range("D6") '169.77'
Set-TextValue $ws.Range("E6") '  -1.96%  '
Set-TextValue $ws.Range("E7") '  +0.05%  '
Set-TextValue $ws.Range("D8") '3.071.45'
Set-TextValue $ws.Range("E8") '  -1.59%  '
Set-TextValue $ws.Range("E9") '  -2.48%  '
Set-TextValue $ws.Range("D10") '6.31'
Set-TextValue $ws.Range("E10") '  -2.16%  '
Set-TextValue $ws.Range("E11") '  -3.14%  '
Set-TextValue $ws.Range("E12") '  -2.66%  '
Set-TextValue $ws.Range("D13") '0.0000238'
Set-TextValue $ws.Range("E13") '  -3.64%  '
Set-TextValue $ws.Range("E14") '  -3.85%  '
Set-TextValue $ws.Range("E15") '  -1.47%  '
Set-TextValue $ws.Range("D16") '3.588.59'
Set-TextValue $ws.Range("E16") '  -1.50%  '
Set-TextValue $ws.Range("D17") '66.166.82'
Set-TextValue $ws.Range("E17") '  -1.19%  '
Set-TextValue $ws.Range("E18") '  -3.00%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D19") '3.079.33'
Set-TextValue $ws.Range("E19") '  -1.44%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D20") '16.56'
Set-TextValue $ws.Range("E20") '  +1.46%  '
Set-TextValue $ws.Range("D21") '484.84'
Set-TextValue $ws.Range("E21") '  +1.71%  '
Set-TextValue $ws.Range("D22") '0.685'
Set-TextValue $ws.Range("E22") '  -3.67%  '
Set-TextValue $ws.Range("D23") '7.66'
Set-TextValue $ws.Range("E23") '  -3.41%  '
Set-TextValue $ws.Range("D24") '82.43'
Set-TextValue $ws.Range("E24") '  -1.85%  '
Set-TextValue $ws.Range("D25") '12.61'
Set-TextValue $ws.Range("E25") '  -5.30%  '
Set-TextValue $ws.Range("E26") '  -4.02%  '
Set-TextValue $ws.Range("D27") '10.22'
Set-TextValue $ws.Range("E27") '  +1.16%  '
Set-TextValue $ws.Range("E28") '  +0.06%  '
Set-TextValue $ws.Range("E29") '  -0.53%  '
Set-TextValue $ws.Range("E30") '  -5.44%  '
Set-TextValue $ws.Range("D31") '2.59'
Set-TextValue $ws.Range("E31") '  -3.19%  '
Set-TextValue $ws.Range("D32") '27.71'
Set-TextValue $ws.Range("E32") '  -3.17%  '
Set-TextValue $ws.Range("D33") '0.110'
Set-TextValue $ws.Range("E33") '  -4.15%  '
Set-TextValue $ws.Range("D34") '0.0₃0897'
Set-TextValue $ws.Range("E34") '  -5.93%  '
Set-TextValue $ws.Range("E35") '  +0.04%  '
Set-TextValue $ws.Range("D36") '0.945'
Set-TextValue $ws.Range("E36") '  -3.37%  '
Set-TextValue $ws.Range("D37") '47.12'
Set-TextValue $ws.Range("E37") '  -0.07%  '
Set-TextValue $ws.Range("D38") '5.55'
Set-TextValue $ws.Range("E38") '  -5.38%  '
Set-TextValue $ws.Range("E39") '  -1.23%  '
Set-TextValue $ws.Range("D40") '1.95'
Set-TextValue $ws.Range("E40") '  -5.37%  '
Set-TextValue $ws.Range("D41") '0.298'
Set-TextValue $ws.Range("E41") '  -4.59%  '
Set-TextValue $ws.Range("D42") '8.24'
Set-TextValue $ws.Range("E42") '  -4.98%  '
Set-TextValue $ws.Range("D43") '2.778.90'
Set-TextValue $ws.Range("E43") '  -1.23%  '
Set-TextValue $ws.Range("D44") '0.0343'
Set-TextValue $ws.Range("E44") '  -3.20%  '
Set-TextValue $ws.Range("E45") '  -2.31%  '
Set-TextValue $ws.Range("D46") '134.72'
Set-TextValue $ws.Range("E46") '  -0.64%  '
Set-TextValue $ws.Range("D47") '363.70'
Set-TextValue $ws.Range("E47") '  -5.23%  '
Set-TextValue $ws.Range("D49") '24.34'
Set-TextValue $ws.Range("E49") '  -2.42%  '
Set-TextValue $ws.Range("D50") '2.14'
Set-TextValue $ws.Range("E50") '  -2.48%  '
Set-TextValue $ws.Range("E51") '  -2.55%  '
